$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos table update - rows 2 to 51, columns D (Price) and E (Volume %)
$ws.Range('D2').Value = '64.348.39'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').Value = '3.505.32'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '586.92'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '134.56'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.488'
$ws.Range('E8').Value = '  -0.27%  '
$ws.Range('E9').Value = '  +0.93%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.23'
$ws.Range('E10').Value = '  +0.31%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.390'
$ws.Range('E11').Value = '  +1.89%  '
$ws.Range('D12').Value = '4.106.40'
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('E13').Value = '  +1.09%  '
$ws.Range('E14').Value = '  +1.35%  '
$ws.Range('D15').Value = '3.507.18'
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '25.83'
$ws.Range('E16').Value = '  -6.59%  '
$ws.Range('D17').Value = '64.346.06'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '9.95'
$ws.Range('E18').Value = '  -3.38%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.79'
$ws.Range('E19').Value = '  +1.58%  '
$ws.Range('E20').Value = '  -4.09%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '394.88'
$ws.Range('E21').Value = '  +1.98%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.572'
$ws.Range('E22').Value = '  -1.16%  '
$ws.Range('D23').Value = '3.646.23'
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '74.39'
$ws.Range('E24').Value = '  +1.84%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  +1.94%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000115'
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.47'
$ws.Range('E28').Value = '  -0.38%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('E31').Value = '  -1.06%  '
$ws.Range('E32').Value = '  -5.44%  '
$ws.Range('D33').Value = '3.527.86'
$ws.Range('E33').Value = '  +0.58%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('E35').Value = '  +2.81%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '23.60'
$ws.Range('E36').Value = '  -0.83%  '
$ws.Range('E37').Value = '  -1.76%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.58'
$ws.Range('E38').Value = '  -0.67%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.92'
$ws.Range('E39').Value = '  -0.44%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '163.01'
$ws.Range('E40').Value = '  -0.53%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0784'
$ws.Range('E41').Value = '  -2.41%  '
$ws.Range('E42').Value = '  -1.21%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '25.86'
$ws.Range('E43').Value = '  -1.70%  '
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '41.86'
$ws.Range('E45').Value = '  +0.90%  '
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('E47').Value = '  -4.09%  '
$ws.Range('E48').Value = '  +0.79%  '
$ws.Range('D49').Value = '2.473.79'
$ws.Range('E49').Value = '  +1.38%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.79'
$ws.Range('E50').Value = '  -1.42%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.896'
$ws.Range('E51').Value = '  +0.10%  '
